$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 174910
$ws.Range("C4").Value = 164896
$ws.Range("C5").Value = 10014
$ws.Range("C7").Value = 5.73
$ws.Range("C8").Value = 64.51000000000001
